# GCSE 9-1 Computer Science course
# Add a new "Sorting and Searching Algorithms" teaching unit: fill in the
# objectives on the existing section-title slide, then add six new slides
# (Bubble sort, Merge sort, Linear search, Binary search, a blank spacer,
# and the next section title "4. Decomposition and Abstraction").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Existing slide 15 ("3. Sorting and Searching Algorithms") - merge the
#    title into a single run and populate the (previously empty) content
#    placeholder with the lesson objectives.
# ---------------------------------------------------------------------
$s15 = $p.Slides.Item(15)

$titleTr = $s15.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "x"
$titleTr.Text = "3. Sorting and Searching Algorithms"

$bodyTr = $s15.Shapes.Item(2).TextFrame.TextRange
$bodyTr.LanguageID = "en-GB"
$bodyTr.Text = "Understand how bubble sort, merge sort, linear search and binary search algorithms work`rUnderstand how the choice of algorithm is influenced by the data structures and data values that need to be manipulated`rEvaluate the fitness for purpose of algorithms in meeting specified requirements efficiently, using logical reasoning and test data.  "

# ---------------------------------------------------------------------
# 2. New slide 16 - "Bubble sort"
# ---------------------------------------------------------------------
$s16 = $p.Slides.Add(16, 2)

$s16.Shapes.Item(1).TextFrame.TextRange.LanguageID = "en-GB"
$s16.Shapes.Item(1).TextFrame.TextRange.Text = "Bubble sort"
$s16.Shapes.Item(1).TextFrame.AutoSize = 2

$c16 = $s16.Shapes.Item(2).TextFrame.TextRange
$c16.LanguageID = "en-GB"
$c16.Text = "Start at the beginning of the list.`rCompare the values in position 1 and position 2 in the list " + [char]0x2013 + " if they are not in ascending (descending) order then swap them. `rCompare the values in position 2 and position 3 in the list and swap if necessary.`rContinue to the end of the list.`rIf there have been any swaps, repeat steps 1 to 4.`rHomework : Exercise : Bubble sort flowchart"

# ---------------------------------------------------------------------
# 3. New slide 17 - "Merge sort"
# ---------------------------------------------------------------------
$s17 = $p.Slides.Add(17, 2)

$s17.Shapes.Item(1).TextFrame.TextRange.LanguageID = "en-GB"
$s17.Shapes.Item(1).TextFrame.TextRange.Text = "Merge sort"

$c17 = $s17.Shapes.Item(2).TextFrame.TextRange
$c17.LanguageID = "en-GB"
$c17.Text = "Divide a list into smaller lists`rDivide these until the size of each list is on`rRecursion is the method to the previous application of the method. `rThe difference between bubble and merge sort is the fact that bubble sort uses brute force, whereas merge sort uses divide and conquer methods to sort the items in a list."

# ---------------------------------------------------------------------
# 4. New slide 18 - "Linear search"
# ---------------------------------------------------------------------
$s18 = $p.Slides.Add(18, 2)

$s18.Shapes.Item(1).TextFrame.TextRange.LanguageID = "en-GB"
$s18.Shapes.Item(1).TextFrame.TextRange.Text = "Linear search"

$c18 = $s18.Shapes.Item(2).TextFrame.TextRange
$c18.LanguageID = "en-GB"
$c18.Text = "Start at the first item in the list. `rCompare the item with the search item.`rIf they are the same, then stop.`rIf they are not, then move to the next item. `rRepeat 2 to 4 until the end of the list is reached. `r"

# ---------------------------------------------------------------------
# 5. New slide 19 - "Binary search"
# ---------------------------------------------------------------------
$s19 = $p.Slides.Add(19, 2)

$s19.Shapes.Item(1).TextFrame.TextRange.LanguageID = "en-GB"
$s19.Shapes.Item(1).TextFrame.TextRange.Text = "Binary search"

$c19 = $s19.Shapes.Item(2).TextFrame.TextRange
$c19.LanguageID = "en-GB"
$c19.Text = "Select the median item of the list. `rIf the median is equal to the search item, then stop.`rIf the media is too high, then repeat 1 and 2 with the sub-list to the left. `rIf the median is too low, then repeat 1 and 2 with the sub-list to the right.`rRepeat steps 3 and 4 until the item has been found or all of the items have been checked.`rBinary search is much more efficient than the linear search."

# ---------------------------------------------------------------------
# 6. New slide 20 - blank spacer slide
# ---------------------------------------------------------------------
$s20 = $p.Slides.Add(20, 2)

# ---------------------------------------------------------------------
# 7. New slide 21 - "4. Decomposition and Abstraction"
# ---------------------------------------------------------------------
$s21 = $p.Slides.Add(21, 2)

$s21.Shapes.Item(1).TextFrame.TextRange.LanguageID = "en-GB"
$s21.Shapes.Item(1).TextFrame.TextRange.Text = "4. Decomposition and Abstraction"
